# LOM3092.xlsx update
#
# The sheet lists course-catalog fields as label/value rows (col A = label,
# cols B/C = value, mirrored). This edit:
#   - inserts a new "Docentes responsáveis:" label row right after
#     "Objectives:" (pushing every following row down by one),
#   - and re-points a handful of value cells that, after the shift, need to
#     show different text than the row they inherited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the old row 12 ("Programa resumido:"); rows
# 12-25 shift down to 13-26.
$ws.Rows("12:12").Insert()

# New row 12 just carries the label in column A.
$ws.Range("A12").Value = "Docentes responsáveis:"

# "Objetivos:" row (10): value becomes the responsible professor line.
$ws.Range("B10").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("C10").Value = "2166002 - Sandra Giacomin Schneider"

# "Programa resumido:" row (13, shifted from 12): value becomes "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "Programa:" row (15, shifted from 14): value becomes the activation date.
# Copy it from B8/C8 (which already holds "01/01/2012" as text) instead of
# typing the literal, so Excel doesn't reinterpret it as a date value/format.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# "Método:" row (18, shifted from 17): value becomes the professor line again.
$ws.Range("B18").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("C18").Value = "2166002 - Sandra Giacomin Schneider"
